# Apply "add user list to project" change:
#  1. Sheet1 "PI hours": Geir Dullerud's cfop list gains an extra code (cfop_DULLERUD_1).
#  2. Sheet4 "project hours": new "users" column (E) listing the user(s) behind each project.
#  3. Sheet5 "cfop hours": new row for cfop_DULLERUD_1 (17 hours), percentages recomputed
#     for the new grand total.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) PI hours: Geir Dullerud (row 4) cfop list now includes the new cfop code.
# ---------------------------------------------------------------------------
$wsPI = $wb.Worksheets.Item("PI hours")
$wsPI.Range("G4").Value = "['cfop_DULLERUD_1', 'cfop_DULLERUD']"

# ---------------------------------------------------------------------------
# 2) project hours: add the "users" column with per-project user lists.
# ---------------------------------------------------------------------------
$wsProj = $wb.Worksheets.Item("project hours")

$wsProj.Range("E1").Value = "users"
# Match the bold/centered/bordered header style already used by B1:D1.
$wsProj.Range("B1").Copy()
$wsProj.Range("E1").PasteSpecial(-4122)   # xlPasteFormats

$projectUsers = @(
    "['Jiyang Chen']",
    "['Alex Hill', 'Alexander Hill']",
    "['Joao Porto']",
    "['Arun Lakshmanan', 'Gabriel Barsi Haberfeld']",
    "['Harshal Maske']",
    "['Usman Syed']",
    "['Karun Koppul']",
    "['Dalton Chaffee']",
    "['Daniel Olivas Hernandez']"
)

for ($i = 0; $i -lt $projectUsers.Length; $i++) {
    $row = $i + 2
    $wsProj.Cells.Item($row, 5).Value = $projectUsers[$i]
}

# ---------------------------------------------------------------------------
# 3) cfop hours: insert the new cfop_DULLERUD_1 row and recompute percentages.
# ---------------------------------------------------------------------------
$wsCfop = $wb.Worksheets.Item("cfop hours")

# Shift the existing rows 6-9 (cfop_NH ... cfop_WORK) down by one row to make
# room for the newly inserted cfop_DULLERUD_1 row.
$wsCfop.Range("A10").Value = 8
$wsCfop.Range("B10").Value = "cfop_WORK"
$wsCfop.Range("C10").Value = 4
$wsCfop.Range("D10").Value = 2.461538461538462

$wsCfop.Range("A9").Value = 7
$wsCfop.Range("B9").Value = "cfop_MITRA"
$wsCfop.Range("C9").Value = 6
$wsCfop.Range("D9").Value = 3.692307692307693

$wsCfop.Range("A8").Value = 6
$wsCfop.Range("B8").Value = "cfop_HUTCHINSON"
$wsCfop.Range("C8").Value = 12
$wsCfop.Range("D8").Value = 7.384615384615385

$wsCfop.Range("A7").Value = 5
$wsCfop.Range("B7").Value = "cfop_NH"
$wsCfop.Range("C7").Value = 16
$wsCfop.Range("D7").Value = 9.846153846153847

# New row for cfop_DULLERUD_1
$wsCfop.Range("A6").Value = 4
$wsCfop.Range("B6").Value = "cfop_DULLERUD_1"
$wsCfop.Range("C6").Value = 17
$wsCfop.Range("D6").Value = 10.46153846153846

# The row-10 shift above left the newly written A10 cell without the
# bold/centered/bordered "index" style used by the rest of column A; restore it.
$wsCfop.Range("A9").Copy()
$wsCfop.Range("A10").PasteSpecial(-4122)   # xlPasteFormats

# Recompute percentages for the rows that stayed in place (totals changed
# from 145.5 to 162.5 hours).
$wsCfop.Range("D2").Value = 25.23076923076923
$wsCfop.Range("D3").Value = 16.92307692307692
$wsCfop.Range("D4").Value = 13.53846153846154
$wsCfop.Range("D5").Value = 10.46153846153846
